$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(28, 8).Value = 2892.6667
$ws.Cells.Item(28, 9).Value = 2289
$ws.Cells.Item(28, 10).Value = 5005.5
$ws.Cells.Item(28, 11).Value = 2289
$ws.Cells.Item(28, 12).Value = 5005.5
$ws.Cells.Item(28, 13).Value = -1804
$ws.Cells.Item(28, 14).Value = -5975.5

$ws.Cells.Item(33, 8).Value = 179.125
$ws.Cells.Item(33, 9).Value = 140.42857
$ws.Cells.Item(33, 11).Value = 140.42857
$ws.Cells.Item(33, 13).Value = 88.57142999999999

$ws.Cells.Item(62, 8).Value = 1098
$ws.Cells.Item(62, 9).Value = 996
$ws.Cells.Item(62, 10).Value = 1200
$ws.Cells.Item(62, 11).Value = 996
$ws.Cells.Item(62, 12).Value = 1200
$ws.Cells.Item(62, 13).Value = -372
$ws.Cells.Item(62, 14).Value = -2448

$ws.Cells.Item(65, 8).Value = 1098
$ws.Cells.Item(65, 9).Value = 996
$ws.Cells.Item(65, 10).Value = 1200
$ws.Cells.Item(65, 11).Value = 4980
$ws.Cells.Item(65, 12).Value = 1200
$ws.Cells.Item(65, 13).Value = -1860
$ws.Cells.Item(65, 14).Value = -12240

$ws.Cells.Item(97, 8).Value = 1200
$ws.Cells.Item(97, 10).Value = 1200
$ws.Cells.Item(97, 12).Value = 3600
$ws.Cells.Item(97, 14).Value = -4592

$ws.Cells.Item(132, 8).Value = 1728.3334
$ws.Cells.Item(132, 9).Value = 1592.5
$ws.Cells.Item(132, 10).Value = 2000
$ws.Cells.Item(132, 11).Value = 4777.5
$ws.Cells.Item(132, 12).Value = 6000
$ws.Cells.Item(132, 13).Value = -2247.5
$ws.Cells.Item(132, 14).Value = -11060

$ws.Cells.Item(135, 8).Value = 1084.9333
$ws.Cells.Item(135, 9).Value = 944.1539
$ws.Cells.Item(135, 11).Value = 8497.3851
$ws.Cells.Item(135, 13).Value = -5962.3851

$ws.Cells.Item(137, 8).Value = 6013.222
$ws.Cells.Item(137, 9).Value = 3032.25
$ws.Cells.Item(137, 10).Value = 8398
$ws.Cells.Item(137, 11).Value = 9096.75
$ws.Cells.Item(137, 12).Value = 25194
$ws.Cells.Item(137, 13).Value = -6546.75
$ws.Cells.Item(137, 14).Value = -30294

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(10, 8).Value = 1838.75
$ws.Cells.Item(10, 10).Value = 1838.75
$ws.Cells.Item(10, 12).Value = 1838.75
$ws.Cells.Item(10, 14).Value = -2178.75

$ws.Cells.Item(32, 8).Value = 571.375
$ws.Cells.Item(32, 9).Value = 571.375
$ws.Cells.Item(32, 11).Value = 571.375
$ws.Cells.Item(32, 13).Value = -284.375

$ws.Cells.Item(38, 8).Value = 994804
$ws.Cells.Item(38, 9).Value = 3010
$ws.Cells.Item(38, 10).Value = 1656000
$ws.Cells.Item(38, 11).Value = 3010
$ws.Cells.Item(38, 12).Value = 1656000
$ws.Cells.Item(38, 13).Value = -2543
$ws.Cells.Item(38, 14).Value = -1656934

$ws.Cells.Item(61, 8).Value = 4399
$ws.Cells.Item(61, 9).Value = 4399
$ws.Cells.Item(61, 11).Value = 4399
$ws.Cells.Item(61, 13).Value = -4187

$ws.Cells.Item(132, 8).Value = 2885
$ws.Cells.Item(132, 9).Value = 0
$ws.Cells.Item(132, 10).Value = 2885
$ws.Cells.Item(132, 11).Value = 0
$ws.Cells.Item(132, 12).Value = $null
$ws.Cells.Item(132, 13).Value = 8655
$ws.Cells.Item(132, 14).Value = -13715

$ws.Cells.Item(136, 8).Value = 4399
$ws.Cells.Item(136, 9).Value = 4399
$ws.Cells.Item(136, 11).Value = 13197
$ws.Cells.Item(136, 13).Value = -10647

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 986
$ws.Cells.Item(20, 9).Value = 650
$ws.Cells.Item(20, 10).Value = 1210
$ws.Cells.Item(20, 11).Value = 650
$ws.Cells.Item(20, 12).Value = 1210
$ws.Cells.Item(20, 13).Value = -403
$ws.Cells.Item(20, 14).Value = -1704

$ws.Cells.Item(86, 8).Value = 2300
$ws.Cells.Item(86, 9).Value = 2300
$ws.Cells.Item(86, 10).Value = 0
$ws.Cells.Item(86, 11).Value = 2300
$ws.Cells.Item(86, 12).Value = 0
$ws.Cells.Item(86, 13).Value = $null
$ws.Cells.Item(86, 14).Value = -1177

$ws.Cells.Item(89, 8).Value = 2300
$ws.Cells.Item(89, 9).Value = 2300
$ws.Cells.Item(89, 10).Value = 0
$ws.Cells.Item(89, 11).Value = 11500
$ws.Cells.Item(89, 12).Value = 0
$ws.Cells.Item(89, 13).Value = $null
$ws.Cells.Item(89, 14).Value = -5884

$ws.Cells.Item(134, 8).Value = 2331.8
$ws.Cells.Item(134, 9).Value = 2331.8
$ws.Cells.Item(134, 10).Value = 0
$ws.Cells.Item(134, 11).Value = 6995.400000000001
$ws.Cells.Item(134, 12).Value = 0
$ws.Cells.Item(134, 13).Value = $null
$ws.Cells.Item(134, 14).Value = -4460.400000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(12, 8).Value = 5
$ws.Cells.Item(12, 9).Value = 5
$ws.Cells.Item(12, 10).Value = 0
$ws.Cells.Item(12, 11).Value = 5
$ws.Cells.Item(12, 12).Value = 0
$ws.Cells.Item(12, 13).Value = $null
$ws.Cells.Item(12, 14).Value = 165

$ws.Cells.Item(122, 8).Value = 3559.111
$ws.Cells.Item(122, 9).Value = 835.3333
$ws.Cells.Item(122, 10).Value = 9006.666999999999
$ws.Cells.Item(122, 11).Value = 2505.9999
$ws.Cells.Item(122, 12).Value = 27020.001
$ws.Cells.Item(122, 13).Value = -55.9998999999998
$ws.Cells.Item(122, 14).Value = -31920.001

$ws.Cells.Item(132, 8).Value = 2979
$ws.Cells.Item(132, 9).Value = 2979
$ws.Cells.Item(132, 11).Value = 8937
$ws.Cells.Item(132, 13).Value = -6407

$ws.Cells.Item(134, 8).Value = 7487.143
$ws.Cells.Item(134, 9).Value = 2082
$ws.Cells.Item(134, 11).Value = 6246
$ws.Cells.Item(134, 13).Value = -3711

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(64, 8).Value = 6006
$ws.Cells.Item(64, 9).Value = 6006
$ws.Cells.Item(64, 10).Value = 0
$ws.Cells.Item(64, 11).Value = 18018
$ws.Cells.Item(64, 12).Value = 0
$ws.Cells.Item(64, 13).Value = $null
$ws.Cells.Item(64, 14).Value = -17748

$ws.Cells.Item(67, 8).Value = 6006
$ws.Cells.Item(67, 9).Value = 6006
$ws.Cells.Item(67, 10).Value = 0
$ws.Cells.Item(67, 11).Value = 18018
$ws.Cells.Item(67, 12).Value = 0
$ws.Cells.Item(67, 13).Value = $null
$ws.Cells.Item(67, 14).Value = -17082

$ws.Cells.Item(137, 8).Value = 0
$ws.Cells.Item(137, 9).Value = 0
$ws.Cells.Item(137, 10).Value = 0
$ws.Cells.Item(137, 11).Value = 0
$ws.Cells.Item(137, 12).Value = $null
$ws.Cells.Item(137, 13).Value = $null
$ws.Cells.Item(137, 14).Value = 0

$ws.Cells.Item(140, 8).Value = 957.4286
$ws.Cells.Item(140, 9).Value = 957.4286
$ws.Cells.Item(140, 11).Value = 2872.2858
$ws.Cells.Item(140, 13).Value = 2307.7142

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(34, 8).Value = 12500
$ws.Cells.Item(34, 10).Value = 12500
$ws.Cells.Item(34, 12).Value = 12500
$ws.Cells.Item(34, 14).Value = -13036

$ws.Cells.Item(76, 8).Value = 12500
$ws.Cells.Item(76, 10).Value = 12500
$ws.Cells.Item(76, 12).Value = 12500
$ws.Cells.Item(76, 14).Value = -13130

$ws.Cells.Item(79, 8).Value = 12500
$ws.Cells.Item(79, 10).Value = 12500
$ws.Cells.Item(79, 12).Value = 12500
$ws.Cells.Item(79, 14).Value = -14684

$ws.Cells.Item(97, 8).Value = 501110.84
$ws.Cells.Item(97, 9).Value = 416.25
$ws.Cells.Item(97, 11).Value = 416.25
$ws.Cells.Item(97, 13).Value = 79.75

$ws.Cells.Item(122, 8).Value = 1986.7778
$ws.Cells.Item(122, 9).Value = 1184.1333
$ws.Cells.Item(122, 10).Value = 6000
$ws.Cells.Item(122, 11).Value = 3552.3999
$ws.Cells.Item(122, 12).Value = 18000
$ws.Cells.Item(122, 13).Value = -1102.3999
$ws.Cells.Item(122, 14).Value = -22900

$ws.Cells.Item(132, 9).Value = 0
$ws.Cells.Item(132, 11).Value = 0
$ws.Cells.Item(132, 13).Value = $null

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 2674.5
$ws.Cells.Item(7, 9).Value = 2674.5
$ws.Cells.Item(7, 10).Value = 0
$ws.Cells.Item(7, 11).Value = 2674.5
$ws.Cells.Item(7, 12).Value = 0
$ws.Cells.Item(7, 13).Value = $null
$ws.Cells.Item(7, 14).Value = -2562.5

$ws.Cells.Item(61, 8).Value = 436.2
$ws.Cells.Item(61, 9).Value = 436.2
$ws.Cells.Item(61, 11).Value = 436.2
$ws.Cells.Item(61, 13).Value = -234.2

$ws.Cells.Item(70, 8).Value = 52000
$ws.Cells.Item(70, 10).Value = 52000
$ws.Cells.Item(70, 12).Value = 52000
$ws.Cells.Item(70, 14).Value = -52540

$ws.Cells.Item(73, 8).Value = 52000
$ws.Cells.Item(73, 10).Value = 52000
$ws.Cells.Item(73, 12).Value = 52000
$ws.Cells.Item(73, 14).Value = -53872

$ws.Cells.Item(113, 8).Value = 436.2
$ws.Cells.Item(113, 9).Value = 436.2
$ws.Cells.Item(113, 11).Value = 436.2
$ws.Cells.Item(113, 13).Value = 1733.8

$ws.Cells.Item(126, 8).Value = 2674.5
$ws.Cells.Item(126, 9).Value = 2674.5
$ws.Cells.Item(126, 10).Value = 0
$ws.Cells.Item(126, 11).Value = 8023.5
$ws.Cells.Item(126, 12).Value = 0
$ws.Cells.Item(126, 13).Value = $null
$ws.Cells.Item(126, 14).Value = -5553.5

$ws.Cells.Item(132, 8).Value = 1624.25
$ws.Cells.Item(132, 9).Value = 1624.25
$ws.Cells.Item(132, 11).Value = 4872.75
$ws.Cells.Item(132, 13).Value = -2342.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(10, 8).Value = 0
$ws.Cells.Item(10, 9).Value = 0
$ws.Cells.Item(10, 10).Value = 0
$ws.Cells.Item(10, 11).Value = 0
$ws.Cells.Item(10, 12).Value = $null
$ws.Cells.Item(10, 13).Value = $null
$ws.Cells.Item(10, 14).Value = 0

$ws.Cells.Item(32, 8).Value = 10026
$ws.Cells.Item(32, 9).Value = 10026
$ws.Cells.Item(32, 11).Value = 10026
$ws.Cells.Item(32, 13).Value = -9709

$ws.Cells.Item(122, 8).Value = 1154.4333
$ws.Cells.Item(122, 9).Value = 1129.75
$ws.Cells.Item(122, 11).Value = 3389.25
$ws.Cells.Item(122, 13).Value = -939.25

$ws.Cells.Item(126, 8).Value = 2131.9
$ws.Cells.Item(126, 9).Value = 2131.9
$ws.Cells.Item(126, 11).Value = 6395.700000000001
$ws.Cells.Item(126, 13).Value = -3925.700000000001

$ws.Cells.Item(132, 8).Value = 712
$ws.Cells.Item(132, 9).Value = 640
$ws.Cells.Item(132, 11).Value = 1920
$ws.Cells.Item(132, 13).Value = 610

$ws.Cells.Item(136, 8).Value = 1391.9678
$ws.Cells.Item(136, 9).Value = 1298.8572
$ws.Cells.Item(136, 10).Value = 1587.5
$ws.Cells.Item(136, 11).Value = 3896.5716
$ws.Cells.Item(136, 12).Value = 4762.5
$ws.Cells.Item(136, 13).Value = -1346.5716
$ws.Cells.Item(136, 14).Value = -9862.5
